# Update "paises.xlsx" country stats + timestamp, per commit
# "Update countries & provincias Spain".
#
# The source feed was refreshed: most rows simply get new totals, but the
# table is kept sorted descending by "Casos totales" (column B), so two
# pairs of neighbouring countries (Etiopia/Costa Rica and Mali/Malta)
# swapped rank - the higher row keeps the new, bigger numbers and the
# country that drops a place keeps the values it had before this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 20:43"

# --- Straight numeric refreshes (country stays on the same row) -------

# Estados Unidos
$ws.Cells.Item(4, 2).Value = 7520529
$ws.Cells.Item(4, 3).Value = 25858
$ws.Cells.Item(4, 4).Value = 4755238
$ws.Cells.Item(4, 5).Value = 2552136
$ws.Cells.Item(4, 7).Value = 495
$ws.Cells.Item(4, 8).Value = 213155

# Brasil
$ws.Cells.Item(6, 2).Value = 4856024
$ws.Cells.Item(6, 3).Value = 6795
$ws.Cells.Item(6, 5).Value = 498286
$ws.Cells.Item(6, 7).Value = 199
$ws.Cells.Item(6, 8).Value = 144966

# España
$ws.Cells.Item(10, 2).Value = 810807
$ws.Cells.Item(10, 3).Value = 3722
$ws.Cells.Item(10, 7).Value = 113
$ws.Cells.Item(10, 8).Value = 32086

# Francia
$ws.Cells.Item(14, 2).Value = 589653
$ws.Cells.Item(14, 3).Value = 12148
$ws.Cells.Item(14, 4).Value = 97778
$ws.Cells.Item(14, 5).Value = 459720
$ws.Cells.Item(14, 7).Value = 136
$ws.Cells.Item(14, 8).Value = 32155

# Alemania
$ws.Cells.Item(25, 2).Value = 297864
$ws.Cells.Item(25, 3).Value = 2334
$ws.Cells.Item(25, 5).Value = 28772
$ws.Cells.Item(25, 7).Value = 6
$ws.Cells.Item(25, 8).Value = 9592

# Argelia
$ws.Cells.Item(63, 2).Value = 51847
$ws.Cells.Item(63, 3).Value = 157
$ws.Cells.Item(63, 4).Value = 36385
$ws.Cells.Item(63, 5).Value = 13713
$ws.Cells.Item(63, 7).Value = 8
$ws.Cells.Item(63, 8).Value = 1749

# Irlanda
$ws.Cells.Item(74, 2).Value = 37063
$ws.Cells.Item(74, 3).Value = 466
$ws.Cells.Item(74, 5).Value = 11898
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 1801

# Guayana Francesa
$ws.Cells.Item(106, 2).Value = 9968
$ws.Cells.Item(106, 3).Value = 2
$ws.Cells.Item(106, 4).Value = 9637
$ws.Cells.Item(106, 5).Value = 264

# Gambia
$ws.Cells.Item(140, 2).Value = 3585
$ws.Cells.Item(140, 3).Value = 1
$ws.Cells.Item(140, 4).Value = 2224
$ws.Cells.Item(140, 5).Value = 1246
$ws.Cells.Item(140, 7).Value = 2
$ws.Cells.Item(140, 8).Value = 115

# Sudan del Sur
$ws.Cells.Item(149, 2).Value = 2715
$ws.Cells.Item(149, 3).Value = 11
$ws.Cells.Item(149, 5).Value = 1375
$ws.Cells.Item(149, 7).Value = 1
$ws.Cells.Item(149, 8).Value = 50

# Vietnam
$ws.Cells.Item(168, 2).Value = 1096
$ws.Cells.Item(168, 3).Value = 1
$ws.Cells.Item(168, 5).Value = 41

# Monaco
$ws.Cells.Item(189, 2).Value = 221
$ws.Cells.Item(189, 3).Value = 2
$ws.Cells.Item(189, 4).Value = 188
$ws.Cells.Item(189, 5).Value = 31

# --- Rank swaps: Etiopia overtakes Costa Rica --------------------------
# Row 53 now holds Etiopia with its refreshed totals; row 54 now holds
# Costa Rica, keeping the totals it had before the refresh.
$ws.Cells.Item(53, 1).Value = "Etiopia"
$ws.Cells.Item(53, 2).Value = 76988
$ws.Cells.Item(53, 3).Value = 890
$ws.Cells.Item(53, 4).Value = 31677
$ws.Cells.Item(53, 5).Value = 44103
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = 1208

$ws.Cells.Item(54, 1).Value = "Costa Rica"
$ws.Cells.Item(54, 2).Value = 76828
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 39843
$ws.Cells.Item(54, 5).Value = 36068
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 917

# --- Rank swaps: Mali overtakes Malta -----------------------------------
# Row 145 now holds Mali with its refreshed totals; row 146 now holds
# Malta, keeping the totals it had before the refresh.
$ws.Cells.Item(145, 1).Value = "Mali"
$ws.Cells.Item(145, 2).Value = 3156
$ws.Cells.Item(145, 3).Value = 25
$ws.Cells.Item(145, 4).Value = 2467
$ws.Cells.Item(145, 5).Value = 558
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 131

$ws.Cells.Item(146, 1).Value = "Malta"
$ws.Cells.Item(146, 2).Value = 3139
$ws.Cells.Item(146, 3).Value = 44
$ws.Cells.Item(146, 4).Value = 2668
$ws.Cells.Item(146, 5).Value = 434
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 2
$ws.Cells.Item(146, 8).Value = 37
